$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("F2").Value = 1.28
$ws.Range("H2").Value = 15.5
$ws.Range("I2").Value = 22
$ws.Range("K2").Value = 6
$ws.Range("P2").Value = 1.88
$ws.Range("Q2").Value = 1.98
$ws.Range("H6").Value = 17
$ws.Range("J6").Value = 8.800000000000001
$ws.Range("K6").Value = 9.199999999999999
$ws.Range("N6").Value = 7.6
$ws.Range("O6").Value = 1.13
$ws.Range("P6").Value = 3.2
$ws.Range("Q6").Value = 1.41
$ws.Range("R6").Value = 1.9
$ws.Range("T6").Value = 2.08
$ws.Range("U6").Value = 1.84
$ws.Range("AC6").Value = 980
$ws.Range("AF6").Value = 9.4
$ws.Range("AJ6").Value = 9.199999999999999
$ws.Range("AK6").Value = 13
$ws.Range("AL6").Value = 980
$ws.Range("AM6").Value = 210
$ws.Range("F7").Value = 2.84
$ws.Range("I7").Value = 2.88
$ws.Range("K7").Value = 3.35
$ws.Range("N7").Value = 3.25
$ws.Range("O7").Value = 1.43
$ws.Range("R7").Value = 1.27
$ws.Range("T7").Value = 1.92
$ws.Range("U7").Value = 1.99
$ws.Range("X7").Value = 11
$ws.Range("Y7").Value = 9.800000000000001
$ws.Range("Z7").Value = 18
$ws.Range("AB7").Value = 9.800000000000001
$ws.Range("AC7").Value = 7
$ws.Range("AF7").Value = 18.5
$ws.Range("AH7").Value = 27
$ws.Range("F8").Value = 1.56
$ws.Range("G8").Value = 1.85
$ws.Range("H8").Value = 2.18
$ws.Range("I8").Value = 11
$ws.Range("Q8").Value = 1.01
$ws.Range("F9").Value = 2.16
$ws.Range("H9").Value = 2.72
$ws.Range("K9").Value = 6.2
$ws.Range("P9").Value = 2.14
$ws.Range("Q9").Value = 1.01
$ws.Range("F10").Value = 2.08
$ws.Range("H10").Value = 2.72
$ws.Range("J10").Value = 3.35
$ws.Range("P10").Value = 2.22
$ws.Range("Q10").Value = 1.01
$ws.Range("J12").Value = 4
$ws.Range("F15").Value = 3.2
$ws.Range("G15").Value = 3.3
$ws.Range("H15").Value = 2.26
$ws.Range("I15").Value = 2.28
$ws.Range("J15").Value = 3.85
$ws.Range("K15").Value = 4
$ws.Range("O15").Value = 1.23
$ws.Range("P15").Value = 2.34
$ws.Range("Q15").Value = 1.69
$ws.Range("R15").Value = 1.54
$ws.Range("S15").Value = 2.74
$ws.Range("T15").Value = 1.6
$ws.Range("X15").Value = 21
$ws.Range("AA15").Value = 32
$ws.Range("AC15").Value = 8.800000000000001
$ws.Range("AE15").Value = 1000
$ws.Range("AF15").Value = 38
$ws.Range("AG15").Value = 17
$ws.Range("AH15").Value = 16.5
$ws.Range("AJ15").Value = 1000
$ws.Range("AL15").Value = 42
$ws.Range("AM15").Value = 80
$ws.Range("AO15").Value = 14
$ws.Range("F16").Value = 1.64
$ws.Range("G16").Value = 1.65
$ws.Range("H16").Value = 5.9
$ws.Range("N16").Value = 4.2
$ws.Range("Q16").Value = 1.83
$ws.Range("R16").Value = 1.43
$ws.Range("S16").Value = 3.15
$ws.Range("T16").Value = 1.87
$ws.Range("U16").Value = 2.02
$ws.Range("X16").Value = 18
$ws.Range("AB16").Value = 9.199999999999999
$ws.Range("AD16").Value = 24
$ws.Range("AG16").Value = 9.800000000000001
$ws.Range("AJ16").Value = 20
$ws.Range("AL16").Value = 38
$ws.Range("AO16").Value = 110
$ws.Range("F17").Value = 4.7
$ws.Range("G17").Value = 4.9
$ws.Range("H17").Value = 1.8
$ws.Range("I17").Value = 1.82
$ws.Range("N17").Value = 5.1
$ws.Range("P17").Value = 2.38
$ws.Range("Q17").Value = 1.64
$ws.Range("R17").Value = 1.55
$ws.Range("S17").Value = 2.6
$ws.Range("T17").Value = 1.67
$ws.Range("U17").Value = 2.38
$ws.Range("X17").Value = 32
$ws.Range("Y17").Value = 11.5
$ws.Range("Z17").Value = 13
$ws.Range("AA17").Value = 22
$ws.Range("AB17").Value = 30
$ws.Range("AC17").Value = 10.5
$ws.Range("AE17").Value = 19.5
$ws.Range("AF17").Value = 1000
$ws.Range("AG17").Value = 19.5
$ws.Range("AN17").Value = 1000
$ws.Range("AO17").Value = 9
$ws.Range("F18").Value = 4.1
$ws.Range("G18").Value = 4.4
$ws.Range("H18").Value = 1.97
$ws.Range("I18").Value = 2.02
$ws.Range("J18").Value = 3.8
$ws.Range("N18").Value = 4.2
$ws.Range("O18").Value = 1.28
$ws.Range("Q18").Value = 1.85
$ws.Range("S18").Value = 3.15
$ws.Range("T18").Value = 1.74
$ws.Range("U18").Value = 2.2
$ws.Range("AA18").Value = 30
$ws.Range("AC18").Value = 8.6
$ws.Range("AE18").Value = 25
$ws.Range("AI18").Value = 980
$ws.Range("AN18").Value = 980
$ws.Range("F19").Value = 4.9
$ws.Range("I19").Value = 1.75
$ws.Range("K19").Value = 4.4
$ws.Range("S19").Value = 2.98
$ws.Range("Y19").Value = 9.800000000000001
$ws.Range("Z19").Value = 11.5
$ws.Range("AA19").Value = 19.5
$ws.Range("AC19").Value = 9.6
$ws.Range("AD19").Value = 9.800000000000001
$ws.Range("AE19").Value = 19
$ws.Range("AF19").Value = 50
$ws.Range("AO19").Value = 9.6
$ws.Range("F20").Value = 5.5
$ws.Range("I20").Value = 1.66
$ws.Range("K20").Value = 4.6
$ws.Range("Q20").Value = 1.59
$ws.Range("R20").Value = 1.66
$ws.Range("S20").Value = 2.42
$ws.Range("T20").Value = 1.64
$ws.Range("U20").Value = 2.42
$ws.Range("Y20").Value = 12.5
$ws.Range("Z20").Value = 12.5
$ws.Range("AA20").Value = 23
$ws.Range("AE20").Value = 16
$ws.Range("AF20").Value = 980
$ws.Range("AI20").Value = 27
$ws.Range("AJ20").Value = 150
$ws.Range("AK20").Value = 75
$ws.Range("AM20").Value = 70
$ws.Range("AN20").Value = 980
$ws.Range("AO20").Value = 6.8
$ws.Range("I21").Value = 6.8
$ws.Range("N21").Value = 3.75
$ws.Range("Q21").Value = 2
$ws.Range("S21").Value = 3.6
$ws.Range("T21").Value = 2
$ws.Range("U21").Value = 1.89
$ws.Range("Y21").Value = 980
$ws.Range("AA21").Value = 240
$ws.Range("AB21").Value = 8.6
$ws.Range("AE21").Value = 1000
$ws.Range("AG21").Value = 9.800000000000001
$ws.Range("AI21").Value = 120
$ws.Range("I23").Value = 5.1
$ws.Range("N23").Value = 2.68
$ws.Range("O23").Value = 1.51
$ws.Range("S23").Value = 5.1
$ws.Range("AB23").Value = 970
$ws.Range("AC23").Value = 970
$ws.Range("AJ23").Value = 26
$ws.Range("F24").Value = 2.04
$ws.Range("G24").Value = 2.38
$ws.Range("H24").Value = 4
$ws.Range("I24").Value = 5.2
$ws.Range("J24").Value = 2.78
$ws.Range("N24").Value = 2.72
$ws.Range("P24").Value = 1.58
$ws.Range("Q24").Value = 2.36
$ws.Range("S24").Value = 4.2
$ws.Range("V24").Value = 1.25
$ws.Range("W24").Value = 1.72
$ws.Range("AG24").Value = 13.5
$ws.Range("AK24").Value = 34
